# Updates cryptocurrency price and 1h-volume-change figures on the
# "cryptos" sheet, mirroring the scheduled GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (plain data cell, no explicit formatting) used to
# keep forced-text price cells visually identical to the rest of the sheet.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "29.206.13"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.859.58"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'240.54"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.07753"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "'25.09"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "1.882.66"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "'5.232"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "'0.7160"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "'90.41"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "29.208.57"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "'5.864"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "'244.56"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "'0.000007800"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "2.104.98"
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'7.995"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  +2.96%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'0.1590"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  +2.34%  "
$ws.Range("D26").Value = "'162.17"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  -0.72%  "
$ws.Range("D27").Value = "'8.906"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "'18.35"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "'1.496"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("D30").Value = "'1.320"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  -2.75%  "
$ws.Range("D31").Value = "'4.418"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("D32").Value = "'4.223"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +3.32%  "
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").Value = "'1.913"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("D35").Value = "'1.171"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("D36").Value = "'0.7275"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("D37").Value = "'2.674"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "'0.01853"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").Value = "'2.686"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").Value = "1.163.21"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("D41").Value = "'0.9061"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  -0.50%  "
$ws.Range("D42").Value = "'6.156"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("D43").Value = "'72.42"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("D44").Value = "'0.9997"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'101.67"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("D46").Value = "2.003.80"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "'0.5219"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").Value = "'9.301"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +1.70%  "
$ws.Range("D51").Value = "'2.860"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  +1.05%  "
